$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values updated per the crypto-price refresh diff.
# Numeric-looking price values in column D are written with a temporary
# Text number-format (reset to the default "Normal" style right after) so
# Excel keeps storing them as text, exactly like the original cells, instead
# of silently converting them into numeric cells.

$ws.Range("D2").Value = "59.624.62"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "2.290.73"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.64%  "
$ws.Range("D9").Value = "2.289.16"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.20%  "
$ws.Range("D15").Value = "59.570.35"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "2.696.46"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").Value = "2.289.15"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.168"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("E32").Value = "  -1.58%  "
$ws.Range("D33").Value = "0.0₃0720"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.376"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "314.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.70%  "
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "135.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0933"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0486"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("E50").Value = "  +18.55%  "
$ws.Range("E51").Value = "  -0.85%  "
